# Regenerate the localization-status handoff report with the new
# e2e test-doc id (853c4be3-eabb-48e4-a9f7-d0dac4b56d08 -> 69f57c00-3142-43d1-8e83-b210ad9ab90c)
# and refreshed handoff timestamps / xliff file names.

$wb = $excel.ActiveWorkbook

$oldId = "853c4be3-eabb-48e4-a9f7-d0dac4b56d08"
$newId = "69f57c00-3142-43d1-8e83-b210ad9ab90c"

# The hyperlink target (commit-pinned URL) is left untouched by the
# original change -- only the displayed text is refreshed.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80456b1d6043c858edaed89327190ce2af94bfdb/e2e/$oldId.md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, "", "", "e2e\$newId.md")

$wsOverview.Range("G2").Value = "2016-09-02 13:07:20"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkUrl, "", "", "$newId.md")

$wsZhCn.Range("G2").Value = "$newId.8407c903d4f40f0933fdbca2d05fba4e50c0eb26.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-02 13:07:15"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkUrl, "", "", "$newId.md")

$wsDeDe.Range("G2").Value = "$newId.8407c903d4f40f0933fdbca2d05fba4e50c0eb26.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-02 13:07:20"
